$d = $word.ActiveDocument
$wdMain = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Get-ParaRangeByText($text) {
    $r = $d.Content
    $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $text"
    }
    return $r.Paragraphs(1).Range
}

# --- Change 1: new "Adding AMcharts to my project" bullet after "Challenge with Unix dates" ---
$anchor = Get-ParaRangeByText("Challenge with Unix dates")
$anchor.InsertParagraphAfter()
$newPara = $d.Range($anchor.End, $anchor.End).Paragraphs(1).Range
$xmlAmcharts = '<w:p xmlns:w="' + $wdMain + '">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Adding </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>AMcharts</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> to my project</w:t></w:r>' +
    '</w:p>'
$newPara.InsertXML($xmlAmcharts)

# --- Change 2: "Complete database persist..." -> "Eliminate duplicate charts in database" ---
$p1 = Get-ParaRangeByText("Complete database persist and render of chart data (working on serializers)")
$xmlEliminate = '<w:p xmlns:w="' + $wdMain + '">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Eliminate duplicate charts in database</w:t></w:r>' +
    '</w:p>'
$p1.InsertXML($xmlEliminate)

# --- Change 3: "Render chart from database..." -> "X " + "Display company data..." ---
$p2 = Get-ParaRangeByText("Render chart from database data returned above in companycontainer")
$xmlDisplay = '<w:p xmlns:w="' + $wdMain + '">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">X </w:t></w:r>' +
    '<w:r><w:t>Display company data in company container (already added to store)</w:t></w:r>' +
    '</w:p>'
$p2.InsertXML($xmlDisplay)

# --- Change 4: merge "X Display company data..." (old run) + "Fetch and display..." paragraph ---
# After change 3, there are now two consecutive paragraphs each starting with an
# "X "-led display-company-data / fetch-chart message. We need the *second*
# occurrence of "Display company data..." (the original one) merged with the
# "Fetch and display stock chart in company container" paragraph that follows it.
$rSearch = $d.Content
$rSearch.Start = $p2.End
$found = $rSearch.Find.Execute("Display company data in company container (already added to store)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find second Display-company-data paragraph" }
$pX = $rSearch.Paragraphs(1).Range
$pFetch = $d.Range($pX.End, $pX.End).Paragraphs(1).Range
$mergedRange = $d.Range($pX.Start, $pFetch.End)
$xmlFetch = '<w:p xmlns:w="' + $wdMain + '">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">X </w:t></w:r>' +
    '<w:r><w:t>Fetch and display stock chart in company container</w:t></w:r>' +
    '</w:p>'
$mergedRange.InsertXML($xmlFetch)

# --- Change 5: prepend "X " run before "Download chart data and load to database / view" ---
$p5 = Get-ParaRangeByText("Download chart data and load to database / view")
$xmlDownload = '<w:p xmlns:w="' + $wdMain + '">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="8"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">X </w:t></w:r>' +
    '<w:r><w:t>Download chart data and load to database / view</w:t></w:r>' +
    '</w:p>'
$p5.InsertXML($xmlDownload)

Write-Output "done"
